$wb = $excel.ActiveWorkbook

$wsNew = $wb.Worksheets.Item("NewImportLogic")
$wsProd = $wb.Worksheets.Item("prodfix")

# --- prodfix sheet: insert a new column I (Expected_File_names) ---
$wsProd.Columns("I").Insert()

# Try to keep column I the same width as column H (cosmetic, best effort)
$wsProd.Columns("I").ColumnWidth = $wsProd.Columns("H").ColumnWidth

# Column B got wider to fit the new text (closest achievable value; the
# engine stores column widths in whole pixels so exact fractional widths
# cannot always be reproduced)
$wsProd.Columns("B").ColumnWidth = 43

# --- Update the "search date" description text in column B ---
$wsProd.Range("B2").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"
$wsProd.Range("B5").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"
$wsProd.Range("B8").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"
$wsProd.Range("B11").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"

# --- Update header row (row 1) ---
$wsProd.Range("B1").Value = "Population_name"
$wsProd.Range("H1").Value = "Files_to_upload"
$wsProd.Range("I1").Value = "Expected_File_names"

# --- Populate the new "Expected_File_names" column with the file's basename ---
$wsProd.Range("I2").Value = "UtilityOutcome_Feature_Extraction_file_QoL_UtilityData_ECON_NoUtility.xlsx"
$wsProd.Range("I5").Value = "UtilityOutcome_Feature_Extraction_file_ECON_UtilityData_QoL_NoUtility.xlsx"
$wsProd.Range("I8").Value = "UtilityOutcome_Feature_Extraction_file_Both_QoL_ECON_Utility.xlsx"
$wsProd.Range("I11").Value = "UtilityOutcome_Feature_Extraction_file_NegativeScenario_QoL__ECON_NoUtility.xlsx"

# --- Switch the active/selected tab from NewImportLogic to prodfix ---
$wsProd.Activate()
$wsProd.Range("I1:I11").Select()
$wsProd.Application.ActiveWindow.ScrollColumn = 7

$wb.Save()
